# "Add NULL string cells to parser" - set C4 on Sheet1 to the text "NULL"
# (previously an empty/typed cell at C4), then leave the selection on C1
# to match the saved workbook view state.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = "NULL"

$ws.Activate()
$ws.Range("C1").Select()
